$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: extend formatting (styles) from the last existing row (106) down to the new rows (107-115)
# so column A keeps its bold/border style and column E keeps its datetime number format.
$ws.Range("A106:V106").Copy($ws.Range("A107:V115"))

# Step 2: overwrite the copied values with the actual new match data.

# Row 107
$ws.Cells.Item(107, 1).Value2 = 106
$ws.Cells.Item(107, 2).Value2 = "poland"
$ws.Cells.Item(107, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(107, 4).Value2 = "2023-2024"
$ws.Cells.Item(107, 5).Value2 = 45226.75
$ws.Cells.Item(107, 6).Value2 = "Korona Kielce"
$ws.Cells.Item(107, 7).Value2 = 5
$ws.Cells.Item(107, 8).Value2 = "Puszcza"
$ws.Cells.Item(107, 9).Value2 = 3
$ws.Cells.Item(107, 10).Value2 = 1.85
$ws.Cells.Item(107, 11).Value2 = "22/10/2023 12:43"
$ws.Cells.Item(107, 12).Value2 = 1.91
$ws.Cells.Item(107, 13).Value2 = "27/10/2023 17:52"
$ws.Cells.Item(107, 14).Value2 = 3.63
$ws.Cells.Item(107, 15).Value2 = "22/10/2023 12:43"
$ws.Cells.Item(107, 16).Value2 = 3.56
$ws.Cells.Item(107, 17).Value2 = "27/10/2023 17:52"
$ws.Cells.Item(107, 18).Value2 = 4.34
$ws.Cells.Item(107, 19).Value2 = "22/10/2023 12:43"
$ws.Cells.Item(107, 20).Value2 = 4.34
$ws.Cells.Item(107, 21).Value2 = "27/10/2023 17:52"
$ws.Cells.Item(107, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/korona-kielce-puszcza/SQSmQWkI/"

# Row 108
$ws.Cells.Item(108, 1).Value2 = 107
$ws.Cells.Item(108, 2).Value2 = "poland"
$ws.Cells.Item(108, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(108, 4).Value2 = "2023-2024"
$ws.Cells.Item(108, 5).Value2 = 45226.85416666666
$ws.Cells.Item(108, 6).Value2 = "LKS Lodz"
$ws.Cells.Item(108, 7).Value2 = 0
$ws.Cells.Item(108, 8).Value2 = "Gornik Zabrze"
$ws.Cells.Item(108, 9).Value2 = 5
$ws.Cells.Item(108, 10).Value2 = 2.44
$ws.Cells.Item(108, 11).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(108, 12).Value2 = 2.8
$ws.Cells.Item(108, 13).Value2 = "27/10/2023 20:27"
$ws.Cells.Item(108, 14).Value2 = 3.3
$ws.Cells.Item(108, 15).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(108, 16).Value2 = 3.29
$ws.Cells.Item(108, 17).Value2 = "27/10/2023 20:21"
$ws.Cells.Item(108, 18).Value2 = 2.88
$ws.Cells.Item(108, 19).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(108, 20).Value2 = 2.69
$ws.Cells.Item(108, 21).Value2 = "27/10/2023 20:27"
$ws.Cells.Item(108, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/lks-lodz-gornik-zabrze/Yu1BQ01b/"

# Row 109
$ws.Cells.Item(109, 1).Value2 = 108
$ws.Cells.Item(109, 2).Value2 = "poland"
$ws.Cells.Item(109, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(109, 4).Value2 = "2023-2024"
$ws.Cells.Item(109, 5).Value2 = 45227.52083333334
$ws.Cells.Item(109, 6).Value2 = "Warta Poznan"
$ws.Cells.Item(109, 7).Value2 = 1
$ws.Cells.Item(109, 8).Value2 = "Piast Gliwice"
$ws.Cells.Item(109, 9).Value2 = 1
$ws.Cells.Item(109, 10).Value2 = 3.35
$ws.Cells.Item(109, 11).Value2 = "23/10/2023 19:12"
$ws.Cells.Item(109, 12).Value2 = 4.15
$ws.Cells.Item(109, 13).Value2 = "28/10/2023 12:21"
$ws.Cells.Item(109, 14).Value2 = 3.05
$ws.Cells.Item(109, 15).Value2 = "23/10/2023 19:12"
$ws.Cells.Item(109, 16).Value2 = 2.95
$ws.Cells.Item(109, 17).Value2 = "28/10/2023 12:21"
$ws.Cells.Item(109, 18).Value2 = 2.3
$ws.Cells.Item(109, 19).Value2 = "23/10/2023 19:12"
$ws.Cells.Item(109, 20).Value2 = 2.21
$ws.Cells.Item(109, 21).Value2 = "28/10/2023 12:21"
$ws.Cells.Item(109, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/warta-poznan-piast-gliwice/fk27Rtnh/"

# Row 110
$ws.Cells.Item(110, 1).Value2 = 109
$ws.Cells.Item(110, 2).Value2 = "poland"
$ws.Cells.Item(110, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(110, 4).Value2 = "2023-2024"
$ws.Cells.Item(110, 5).Value2 = 45227.625
$ws.Cells.Item(110, 6).Value2 = "Pogon Szczecin"
$ws.Cells.Item(110, 7).Value2 = 2
$ws.Cells.Item(110, 8).Value2 = "Jagiellonia"
$ws.Cells.Item(110, 9).Value2 = 1
$ws.Cells.Item(110, 10).Value2 = 1.75
$ws.Cells.Item(110, 11).Value2 = "24/10/2023 18:12"
$ws.Cells.Item(110, 12).Value2 = 1.93
$ws.Cells.Item(110, 13).Value2 = "28/10/2023 14:53"
$ws.Cells.Item(110, 14).Value2 = 3.99
$ws.Cells.Item(110, 15).Value2 = "24/10/2023 18:12"
$ws.Cells.Item(110, 16).Value2 = 3.79
$ws.Cells.Item(110, 17).Value2 = "28/10/2023 14:53"
$ws.Cells.Item(110, 18).Value2 = 4.09
$ws.Cells.Item(110, 19).Value2 = "24/10/2023 18:12"
$ws.Cells.Item(110, 20).Value2 = 3.94
$ws.Cells.Item(110, 21).Value2 = "28/10/2023 14:56"
$ws.Cells.Item(110, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/pogon-szczecin-jagiellonia/KYpeOAKU/"

# Row 111
$ws.Cells.Item(111, 1).Value2 = 110
$ws.Cells.Item(111, 2).Value2 = "poland"
$ws.Cells.Item(111, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(111, 4).Value2 = "2023-2024"
$ws.Cells.Item(111, 5).Value2 = 45227.72916666666
$ws.Cells.Item(111, 6).Value2 = "Cracovia"
$ws.Cells.Item(111, 7).Value2 = 1
$ws.Cells.Item(111, 8).Value2 = "Lech Poznan"
$ws.Cells.Item(111, 9).Value2 = 1
$ws.Cells.Item(111, 10).Value2 = 3.55
$ws.Cells.Item(111, 11).Value2 = "24/10/2023 18:12"
$ws.Cells.Item(111, 12).Value2 = 3.52
$ws.Cells.Item(111, 13).Value2 = "28/10/2023 17:26"
$ws.Cells.Item(111, 14).Value2 = 3.37
$ws.Cells.Item(111, 15).Value2 = "24/10/2023 18:12"
$ws.Cells.Item(111, 16).Value2 = 3.27
$ws.Cells.Item(111, 17).Value2 = "28/10/2023 17:26"
$ws.Cells.Item(111, 18).Value2 = 2.15
$ws.Cells.Item(111, 19).Value2 = "24/10/2023 18:12"
$ws.Cells.Item(111, 20).Value2 = 2.25
$ws.Cells.Item(111, 21).Value2 = "28/10/2023 17:26"
$ws.Cells.Item(111, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/cracovia-lech-poznan/WbZvShL5/"

# Row 112
$ws.Cells.Item(112, 1).Value2 = 111
$ws.Cells.Item(112, 2).Value2 = "poland"
$ws.Cells.Item(112, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(112, 4).Value2 = "2023-2024"
$ws.Cells.Item(112, 5).Value2 = 45227.83333333334
$ws.Cells.Item(112, 6).Value2 = "Ruch Chorzow"
$ws.Cells.Item(112, 7).Value2 = 2
$ws.Cells.Item(112, 8).Value2 = "Slask Wroclaw"
$ws.Cells.Item(112, 9).Value2 = 2
$ws.Cells.Item(112, 10).Value2 = 2.29
$ws.Cells.Item(112, 11).Value2 = "22/10/2023 15:12"
$ws.Cells.Item(112, 12).Value2 = 2.75
$ws.Cells.Item(112, 13).Value2 = "28/10/2023 19:42"
$ws.Cells.Item(112, 14).Value2 = 3.27
$ws.Cells.Item(112, 15).Value2 = "22/10/2023 15:12"
$ws.Cells.Item(112, 16).Value2 = 3.16
$ws.Cells.Item(112, 17).Value2 = "28/10/2023 19:42"
$ws.Cells.Item(112, 18).Value2 = 3.33
$ws.Cells.Item(112, 19).Value2 = "22/10/2023 15:12"
$ws.Cells.Item(112, 20).Value2 = 2.83
$ws.Cells.Item(112, 21).Value2 = "28/10/2023 19:42"
$ws.Cells.Item(112, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/ruch-chorzow-slask-wroclaw/QFc3SMWo/"

# Row 113
$ws.Cells.Item(113, 1).Value2 = 112
$ws.Cells.Item(113, 2).Value2 = "poland"
$ws.Cells.Item(113, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(113, 4).Value2 = "2023-2024"
$ws.Cells.Item(113, 5).Value2 = 45228.52083333334
$ws.Cells.Item(113, 6).Value2 = "Zaglebie"
$ws.Cells.Item(113, 7).Value2 = 2
$ws.Cells.Item(113, 8).Value2 = "Radomiak Radom"
$ws.Cells.Item(113, 9).Value2 = 3
$ws.Cells.Item(113, 10).Value2 = 2.36
$ws.Cells.Item(113, 11).Value2 = "22/10/2023 13:43"
$ws.Cells.Item(113, 12).Value2 = 2.2
$ws.Cells.Item(113, 13).Value2 = "29/10/2023 12:28"
$ws.Cells.Item(113, 14).Value2 = 3.35
$ws.Cells.Item(113, 15).Value2 = "22/10/2023 13:43"
$ws.Cells.Item(113, 16).Value2 = 3.46
$ws.Cells.Item(113, 17).Value2 = "29/10/2023 12:28"
$ws.Cells.Item(113, 18).Value2 = 3.1
$ws.Cells.Item(113, 19).Value2 = "22/10/2023 13:43"
$ws.Cells.Item(113, 20).Value2 = 3.44
$ws.Cells.Item(113, 21).Value2 = "29/10/2023 12:28"
$ws.Cells.Item(113, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/zaglebie-radomiak-radom/lzUqRCzC/"

# Row 114
$ws.Cells.Item(114, 1).Value2 = 113
$ws.Cells.Item(114, 2).Value2 = "poland"
$ws.Cells.Item(114, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(114, 4).Value2 = "2023-2024"
$ws.Cells.Item(114, 5).Value2 = 45228.625
$ws.Cells.Item(114, 6).Value2 = "Legia"
$ws.Cells.Item(114, 7).Value2 = 1
$ws.Cells.Item(114, 8).Value2 = "Stal Mielec"
$ws.Cells.Item(114, 9).Value2 = 3
$ws.Cells.Item(114, 10).Value2 = 1.38
$ws.Cells.Item(114, 11).Value2 = "24/10/2023 07:42"
$ws.Cells.Item(114, 12).Value2 = 1.28
$ws.Cells.Item(114, 13).Value2 = "29/10/2023 14:57"
$ws.Cells.Item(114, 14).Value2 = 4.98
$ws.Cells.Item(114, 15).Value2 = "24/10/2023 07:42"
$ws.Cells.Item(114, 16).Value2 = 5.96
$ws.Cells.Item(114, 17).Value2 = "29/10/2023 14:57"
$ws.Cells.Item(114, 18).Value2 = 8.24
$ws.Cells.Item(114, 19).Value2 = "24/10/2023 07:42"
$ws.Cells.Item(114, 20).Value2 = 11.47
$ws.Cells.Item(114, 21).Value2 = "29/10/2023 14:57"
$ws.Cells.Item(114, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/legia-stal-mielec/MmSiPj5O/"

# Row 115
$ws.Cells.Item(115, 1).Value2 = 114
$ws.Cells.Item(115, 2).Value2 = "poland"
$ws.Cells.Item(115, 3).Value2 = "ekstraklasa"
$ws.Cells.Item(115, 4).Value2 = "2023-2024"
$ws.Cells.Item(115, 5).Value2 = 45228.72916666666
$ws.Cells.Item(115, 6).Value2 = "Rakow"
$ws.Cells.Item(115, 7).Value2 = 1
$ws.Cells.Item(115, 8).Value2 = "Widzew Lodz"
$ws.Cells.Item(115, 9).Value2 = 1
$ws.Cells.Item(115, 10).Value2 = 1.44
$ws.Cells.Item(115, 11).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(115, 12).Value2 = 1.44
$ws.Cells.Item(115, 13).Value2 = "29/10/2023 17:22"
$ws.Cells.Item(115, 14).Value2 = 4.69
$ws.Cells.Item(115, 15).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(115, 16).Value2 = 4.72
$ws.Cells.Item(115, 17).Value2 = "29/10/2023 17:22"
$ws.Cells.Item(115, 18).Value2 = 7.04
$ws.Cells.Item(115, 19).Value2 = "22/10/2023 20:15"
$ws.Cells.Item(115, 20).Value2 = 7.51
$ws.Cells.Item(115, 21).Value2 = "29/10/2023 17:22"
$ws.Cells.Item(115, 22).Value2 = "https://www.betexplorer.com/football/poland/ekstraklasa/rakow-czestochowa-widzew-lodz/nPbaT2Hu/"
